$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "duplicate_image_filename" (column E) header is in E1; fill the
# practice rows (2-5) and the main trial rows (6-21) below it with "NA".
$ws.Range("E2:E21").Value = "NA"
